$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the effort/time values
$ws.Range("B4").Value = 7
$ws.Range("B6").Value = 6

# Move the active selection to A2
$ws.Range("A2").Select()
